$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old sub-header row (Hiver/Eté/Année/(m3/s)/(MW)/(GWh) units row).
# This shifts all the data rows up by one.
$ws.Rows.Item(2).Delete()

# Rebuild the header row (row 1) with the new column layout: idx, idx2, Name,
# Date Start, Date End, (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
# E1 previously carried the leftover "(Turbinage et ali)mation" style (9pt Arial);
# the new header label goes back to the sheet's plain default formatting.
$ws.Range("E1").Font.Size = 10
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Give the new numeric/unit header cells (F1:K1) the same font formatting used
# throughout the rest of the sheet.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# Update the selection to match the new layout.
$ws.Range("A2:K2").Select()
